$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW50-FE-LIFTER -> add row 72 ---
$ws1 = $wb.Worksheets.Item(1)
$r = 72
$ws1.Cells.Item($r, 1).Value = 45760.25007583333
$ws1.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item($r, 2).Value = "0x01,0x90"
$ws1.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item($r, 4).Value = "0x01,0x52"
$ws1.Cells.Item($r, 5).Value = "0xe"
$ws1.Cells.Item($r, 6).Value = 400
$ws1.Cells.Item($r, 7).Value = 568631262647114000000000.0
$ws1.Cells.Item($r, 8).Value = 338
$ws1.Cells.Item($r, 9).Value = 14

# --- Sheet 2: ROW50-MID-LIFTER -> add row 74 ---
$ws2 = $wb.Worksheets.Item(2)
$r = 74
$ws2.Cells.Item($r, 1).Value = 45760.21266203704
$ws2.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws2.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item($r, 4).Value = "0x01,0x56"
$ws2.Cells.Item($r, 5).Value = "0x19"
$ws2.Cells.Item($r, 6).Value = 400
$ws2.Cells.Item($r, 7).NumberFormat = "@"
$ws2.Cells.Item($r, 7).Value = "568631262647113771663628"
$ws2.Cells.Item($r, 7).ClearFormats()
$ws2.Cells.Item($r, 8).Value = 342
$ws2.Cells.Item($r, 9).Value = 25

# --- Sheet 3: ROW11-FE-LIFTER -> add row 72 ---
$ws3 = $wb.Worksheets.Item(3)
$r = 72
$ws3.Cells.Item($r, 1).Value = 45760.283908125
$ws3.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item($r, 2).Value = "0x01,0x90"
$ws3.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item($r, 4).Value = "0x01,0x52"
$ws3.Cells.Item($r, 5).Value = "0x14"
$ws3.Cells.Item($r, 6).Value = 400
$ws3.Cells.Item($r, 7).Value = 568631262647114000000000.0
$ws3.Cells.Item($r, 8).Value = 338
$ws3.Cells.Item($r, 9).Value = 20

# --- Sheet 4: ROW11-MID-LIFTER -> add row 72 ---
$ws4 = $wb.Worksheets.Item(4)
$r = 72
$ws4.Cells.Item($r, 1).Value = 45760.40883966435
$ws4.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Cells.Item($r, 2).Value = "0x01,0x90"
$ws4.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item($r, 4).Value = "0x01,0x56"
$ws4.Cells.Item($r, 5).Value = "0x19"
$ws4.Cells.Item($r, 6).Value = 400
$ws4.Cells.Item($r, 7).Value = 568631262647114000000000.0
$ws4.Cells.Item($r, 8).Value = 342
$ws4.Cells.Item($r, 9).Value = 25
